# Update "相談件数" (consultation count) sheet with one more day of data.
#
# Before: last data row is row 98 (2020-05-02), row 99 holds the footnote.
# After : a new data row is inserted as row 99 (2020-05-03 / 43954),
#         pushing the footnote down to row 100.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the footnote row, inheriting the formatting of the
# row above it (date/number styles used by every other data row).
$ws.Rows.Item(99).Insert()

# New day's figures.
$ws.Range("A99").Value = 43954
$ws.Range("B99").Value = 308
$ws.Range("C99").Value = 33036
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 6958

# The print area grew by one row (now covers the newly inserted row as well
# as the footnote row that shifted from 99 to 100).
$wb.Names.Item("相談件数!Print_Area").RefersTo = "=相談件数!`$A`$1:`$E`$101"

# Reflect the new selection on the sheet (bottom-right frozen pane).
$ws.Activate()
$ws.Range("D99").Select()
